$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.103.54"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "1.656.18"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'216.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").Value = "'0.5146"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.13%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").Value = "'0.2632"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.25%  "
$ws.Range("D9").Value = "'0.06256"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("D10").Value = "'20.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.86%  "
$ws.Range("D11").Value = "'0.07719"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").Value = "1.655.74"
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("D13").Value = "'4.427"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").Value = "1.883.73"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "'0.5404"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.15%  "
$ws.Range("E16").Value = "  -3.12%  "
$ws.Range("D17").Value = "'64.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "26.146.26"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D20").Value = "'4.593"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.07%  "
$ws.Range("D21").Value = "'191.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").Value = "'10.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("D23").Value = "'5.995"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.16%  "
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").Value = "'139.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  -4.35%  "
$ws.Range("D27").Value = "'7.195"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("D29").Value = "'1.428"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").Value = "'0.05960"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.60%  "
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("D32").Value = "'3.562"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").Value = "'3.249"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.56%  "
$ws.Range("D34").Value = "'1.597"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.75%  "
$ws.Range("D35").Value = "'0.9608"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.20%  "
$ws.Range("D36").Value = "'2.423"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "'2.771"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").Value = "'0.5661"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.82%  "
$ws.Range("D39").Value = "'0.01585"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.73%  "
$ws.Range("D40").Value = "'5.940"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.58%  "
$ws.Range("D41").Value = "'0.8532"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "1.006.55"
$ws.Range("E43").Value = "  -8.04%  "
$ws.Range("D44").Value = "'100.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "1.798.50"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("E46").Value = "  -4.47%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'56.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.43%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.003"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").Value = "'8.002"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("D50").Value = "'0.05169"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").Value = "'0.4182"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.20%  "
